# CambioTasaCuentaAhorro.xlsx - update entregable 1, 2 y 3
# Update the "usuario apro" row (row 2) with a new run's values, and
# adjust the selected cell / column width to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 value updates (A2, C2, G2, H2 change; B2/D2/E2/F2 stay the same)
$ws.Range("A2").Value = "jtangt"
# C2 holds a numeric-looking id that must stay text (column already uses the
# "number stored as text" quote-prefix style) - use a leading apostrophe so
# it isn't reinterpreted as a number.
$ws.Range("C2").Formula = "'1010826124"
$ws.Range("G2").Value = "AAACT2318450ZF62H"
$ws.Range("H2").Value = "3 jul. 2023, 11:29:15"

# Widen column G to fit the new, longer transaction id
$ws.Columns.Item(7).ColumnWidth = 19.6666666666667

# Move the active selection from E6 to D6
$ws.Range("D6").Select() | Out-Null
